# Actualización automática 2025-09-24 14:25:10
# Inserts a new salesperson row ("VERA TREJO JUAN CARLOS") before
# "VIZUETE GALARZA EDWIN RODRIGO" (row 45) on both the "VENTAS POR GRUPO"
# and "VENTA MENSUAL" sheets, shifting subsequent rows down by one, and
# updates the trailing "X de 45" summary row on "VENTAS POR GRUPO" to
# "X de 46" to reflect the new headcount.

$wb = $excel.ActiveWorkbook

$sheetNames = @("VENTAS POR GRUPO", "VENTA MENSUAL")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $lastCol = $ws.Cells.Item(1, $ws.UsedRange.Columns.Count).Column
    $insertRow = 45

    # Push "VIZUETE GALARZA EDWIN RODRIGO" (and everything below it) down
    # one row, leaving a blank row 45 (with formatting inherited from the
    # row being pushed down) ready for the new salesperson.
    $ws.Rows.Item($insertRow).Insert()

    $ws.Cells.Item($insertRow, 1).Value = "OFICINA-CATAECSA"
    $ws.Cells.Item($insertRow, 2).Value = "VERA TREJO JUAN CARLOS"

    for ($col = 3; $col -le $lastCol; $col++) {
        $ws.Cells.Item($insertRow, $col).Value = 0
    }
}

# Update the trailing "X de 45" -> "X de 46" summary row on
# "VENTAS POR GRUPO" (now at row 48 after the insert above).
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$summaryRow = 48
$lastColGrupo = $wsGrupo.Cells.Item(1, $wsGrupo.UsedRange.Columns.Count).Column
for ($col = 3; $col -le $lastColGrupo; $col++) {
    $cell = $wsGrupo.Cells.Item($summaryRow, $col)
    $text = $cell.Value()
    $cell.Value = $text.Replace("de 45", "de 46")
}
